$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy row 2 (the existing "Panobinostat_1" data row) down into row 3,
# carrying over all formatting, then just change the experiment name.
$ws.Range("A2:L2").Copy() | Out-Null
$ws.Range("A3:L3").PasteSpecial(-4104) | Out-Null   # xlPasteAll
$excel.CutCopyMode = 0

$ws.Range("A3").Value = "Panobinostat_2"

# Update the selected/active cell to A4, matching the edited workbook.
$ws.Range("A4").Select() | Out-Null
